$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" footer timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 21:22"

# Row 8 - Alemania
$ws.Range("B8").Value = 131359
$ws.Range("C8").Value = 1287
$ws.Range("E8").Value = 59865
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 3294

# Row 20 - Austria
$ws.Range("B20").Value = 14226
$ws.Range("C20").Value = 185
$ws.Range("E20").Value = 6209

# Row 26 - Peru
$ws.Range("D26").Value = 2869
$ws.Range("E26").Value = 7204
$ws.Range("F26").Value = 132

# Row 47 - Luxemburgo
$ws.Range("E47").Value = 2740
$ws.Range("H47").Value = 67

# Row 93 - Costa Rica
$ws.Range("B93").Value = 618
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 66
$ws.Range("E93").Value = 549
